# Delete column I (COSMO_Screening_Charge). This shifts column J
# (Max_Absorption_nm) into I, and column K (Max_f_osc) into J,
# matching the target layout: A..H unchanged, I=Max_Absorption_nm,
# J=Max_f_osc.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(9).Delete()
